# Auto-generated PowerShell COM-interop edit script.
# Adds the new "24. 8. 2021" survey wave as column AH on sheet "data"
# and as column AG on sheet "pocetR"; updates the two summary title cells.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsPocet = $wb.Worksheets.Item("pocetR")

# ---------------------------------------------------------------------
# Sheet "data": new column AH (percentages)
# ---------------------------------------------------------------------

# Header cell AH1: copy the header style (bold, border, centered) from AG1,
# then write the new date label.
$wsData.Range("AG1").Copy()
$wsData.Range("AH1").PasteSpecial(-4122)
$wsData.Range("AH1").Value = "24. 8. 2021"

# Numeric values for AH2:AH76 (rows 2..76)
$dataVals = @(0.28, 0.45, 0.27, 0.35, 0.46, 0.19, 0.29, 0.46, 0.25, 0.21, 0.42, 0.37, 0.26, 0.46, 0.28, 0.24, 0.42, 0.34, 0.29, 0.45, 0.26, 0.34, 0.42, 0.24, 0.21, 0.48, 0.31, 0.26, 0.47, 0.27, 0.28, 0.41, 0.31, 0.34, 0.39, 0.27, 0.26, 0.47, 0.27, 0.08, 0.66, 0.26, 0.3, 0.45, 0.25, 0.22, 0.46, 0.32, 0.39, 0.54, 0.07000000000000001, 0.36, 0.37, 0.27, 0.31, 0.46, 0.23, 0.4, 0.44, 0.16, 0.36, 0.4, 0.24, 0.29, 0.42, 0.29, 0.29, 0.48, 0.23, 0.22, 0.53, 0.25, 0.16, 0.43, 0.41)
for ($i = 0; $i -lt $dataVals.Length; $i++) {
    $row = 2 + $i
    $wsData.Cells.Item($row, 34).Value = $dataVals[$i]
}

# Update the footnote/title cell A77 with the new "aktualizace" date
$wsData.Range("A77").Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new column AG (sample sizes)
# ---------------------------------------------------------------------

# Header cell AG1: copy the header style from AF1, then write the new date label.
$wsPocet.Range("AF1").Copy()
$wsPocet.Range("AG1").PasteSpecial(-4122)
$wsPocet.Range("AG1").Value = "24. 8. 2021"

# Numeric values for AG2:AG26 (rows 2..26)
$pocetVals = @(1898, 452, 706, 740, 314, 347, 1237, 934, 964, 985, 442, 220, 251, 41, 150, 90, 17, 261, 523, 235, 367, 332, 217, 339, 408)
for ($i = 0; $i -lt $pocetVals.Length; $i++) {
    $row = 2 + $i
    $wsPocet.Cells.Item($row, 33).Value = $pocetVals[$i]
}

# Row 27: update the footnote title, then add the trailing empty cell AG27
# (mirrors the existing empty B27:AF27 placeholder cells).
$wsPocet.Range("A27").Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"
$wsPocet.Range("AF27").Copy()
$wsPocet.Range("AG27").PasteSpecial(-4122)

